# Rename the "default_none" marker value to "default" across the
# stochastic-structure related sheets, and restore the view/selection
# state (active tab, scroll position, selected ranges) to match the
# author's final save.

$wb = $excel.ActiveWorkbook

# --- 1. Content edit: default_none -> default -------------------------
# rel__sto_struc__sto_scen!B2
$wsStoStruc = $wb.Worksheets.Item("rel__sto_struc__sto_scen")
$wsStoStruc.Range("B2").Value = "default"

# rel_node__stochastic_structure!C2:C156
$wsNodeStoStruc = $wb.Worksheets.Item("rel_node__stochastic_structure")
$wsNodeStoStruc.Range("C2:C156").Value = "default"

# rel_for_node_group_ptdf!E75
$wsNodeGroupPtdf = $wb.Worksheets.Item("rel_for_node_group_ptdf")
$wsNodeGroupPtdf.Range("E75").Value = "default"

# --- 2. View / selection state ------------------------------------------
# obj_connection_ptdf: scroll back up, no longer the selected tab
$wsObjConnPtdf = $wb.Worksheets.Item("obj_connection_ptdf")
$wsObjConnPtdf.Activate()
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
$wsObjConnPtdf.Range("C3:C123").Select()

# rel_for_node_group_ptdf: becomes the selected/active tab
$wsNodeGroupPtdf.Activate()
$excel.ActiveWindow.ScrollRow = 52
$excel.ActiveWindow.ScrollColumn = 1
$wsNodeGroupPtdf.Range("E76").Select()

# rel__sto_struc__sto_scen selection
$wsStoStruc.Activate()
$wsStoStruc.Range("B3").Select()

# rel_node__stochastic_structure selection (extend sqref to C2:C156)
$wsNodeStoStruc.Activate()
$wsNodeStoStruc.Range("C2:C156").Select()

# Final active sheet is rel_for_node_group_ptdf, with bus as first
# visible tab in the sheet-tab scroller.
$wsBus = $wb.Worksheets.Item("bus")
$wsNodeGroupPtdf.Activate()
$excel.ActiveWindow.DisplayedWorksheets = $null
$excel.Windows.Item(1).ScrollWorkbookTabs(1, 1) | Out-Null
